# Generate Report for Handoff
# Updates the "Latest Handoff Datetime" (column D) for the
# d16fe0fa-6a86-4d04-b576-48487efeccd5.md row (row 5) on both the
# "zh-cn" and "de-de" localization status sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-03-10 05:27:07"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-03-10 05:27:19"
